$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E13").Value = 36
$ws.Range("E14").Value = 4
$ws.Range("E15").Value = 118
$ws.Range("E16").Value = 15
$ws.Range("E18").Value = 35
$ws.Range("E19").Value = 4
$ws.Range("E21").Value = 132
$ws.Range("E22").Value = 23
$ws.Range("E24").Value = 30
$ws.Range("E25").Value = 7
$ws.Range("E27").Value = 75
$ws.Range("E28").Value = 10
$ws.Range("E33").Value = 479
$ws.Range("E34").Value = 76
$ws.Range("E36").Value = 652
$ws.Range("E37").Value = 120
$ws.Range("E42").Value = 712
$ws.Range("E43").Value = 230
$ws.Range("E45").Value = 862
$ws.Range("E46").Value = 260
$ws.Range("E48").Value = 860
$ws.Range("E49").Value = 278
$ws.Range("E51").Value = 906
$ws.Range("E52").Value = 335
$ws.Range("E54").Value = 858
$ws.Range("E55").Value = 293
$ws.Range("E57").Value = 839
$ws.Range("E58").Value = 399
$ws.Range("E60").Value = 509
$ws.Range("E61").Value = 203
$ws.Range("E63").Value = 969
$ws.Range("E64").Value = 424
$ws.Range("E66").Value = 941
$ws.Range("E67").Value = 542
$ws.Range("E69").Value = 576
$ws.Range("E70").Value = 256
$ws.Range("E72").Value = 724
$ws.Range("E73").Value = 405
$ws.Range("E75").Value = 779
$ws.Range("E76").Value = 447
$ws.Range("E78").Value = 538
$ws.Range("E79").Value = 471
$ws.Range("E81").Value = 691
$ws.Range("E82").Value = 692
$ws.Range("E84").Value = 551
$ws.Range("E85").Value = 547
$ws.Range("E87").Value = 755
$ws.Range("E88").Value = 420
$ws.Range("E90").Value = 534
$ws.Range("E91").Value = 387
$ws.Range("E96").Value = 537
$ws.Range("E97").Value = 365
$ws.Range("E99").Value = 645
$ws.Range("E100").Value = 498
$ws.Range("E102").Value = 724
$ws.Range("E103").Value = 503
$ws.Range("E105").Value = 289
$ws.Range("E106").Value = 340
$ws.Range("E108").Value = 588
$ws.Range("E109").Value = 811
$ws.Range("E110").Value = 30
